$wb = $excel.ActiveWorkbook

# The status text "Ready for handoff" changed to "In Translation" in every
# sheet that shows it: Overview!E2:F2 (the per-language status columns),
# zh-cn!C2 and de-de!C2 (the "Status" column of each handoff table).
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "In Translation"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "In Translation"

# The "Status" columns narrowed along with the shorter text
# (Overview!E:F, zh-cn!C, de-de!C).
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5
$zhcn.Columns.Item(3).ColumnWidth = 12.5
$dede.Columns.Item(3).ColumnWidth = 12.5
